$d = $word.ActiveDocument

# Locate the list item paragraph that ends with "...and made them scalable"
# (the "Added buttons for art Instagram and embroidery Instagram and made
# them scalable" entry). The new list item belongs right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*and made them scalable*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Collapsed range positioned exactly at the end of the target paragraph
    # (i.e. right after its paragraph mark / right before the next paragraph).
    $insertionPoint = $target.Range.End
    $r = $d.Range($insertionPoint, $insertionPoint)

    # Insert a brand-new list-item paragraph with the same list/paragraph
    # formatting used by the surrounding log entries.
    $newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' +
            '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
            '<w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
            '<w:t>Instagram buttons complete with logo</w:t>' +
        '</w:r>' +
    '</w:p>'

    [void]$r.InsertXML($newParagraphXml)
}
